# Update the "view/visit count" column F on the sheets "展览" and "全部类型"
# to reflect newly generated output (gh-pages rebuild).

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> hashtable of row -> new F value
$updates = @{
    "展览" = @{
        3  = 510
        4  = 1497
        10 = 726
        13 = 318
        14 = 48
        15 = 6335
        16 = 83
        20 = 15180
        23 = 134
        25 = 11001
        27 = 4291
    }
    "全部类型" = @{
        3  = 510
        4  = 1497
        11 = 726
        15 = 318
        16 = 48
        18 = 6335
        19 = 83
        23 = 15180
        26 = 134
        28 = 11001
        30 = 4291
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
